$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A35").Value = "Missed Visitor 2"
$ws.Range("B35").Value = "https://squareonemediauk.files.wordpress.com/2012/04/tom-7957.jpg"

$ws.Range("B35").Select()
